# Weekly fruit/vegetable price update:
# insert two new rows (one new week of data) right before row 435,
# pushing the existing rows 435-459 down to 437-461.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(435).Resize(2).Insert()

# New row 435: Primera quality for the new week (2021-11-16 = serial 44516)
$ws.Cells.Item(435,1).Value = 9
$ws.Cells.Item(435,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(435,3).Value = "Metropolitana"
$ws.Cells.Item(435,4).Value = 44516
$ws.Cells.Item(435,5).Value = 13
$ws.Cells.Item(435,6).Value = 100112023
$ws.Cells.Item(435,7).Value = "Brócoli"
$ws.Cells.Item(435,8).Value = "Sin especificar"
$ws.Cells.Item(435,9).Value = "Primera"
$ws.Cells.Item(435,10).Value = 5200
$ws.Cells.Item(435,11).Value = 500
$ws.Cells.Item(435,12).Value = 600
$ws.Cells.Item(435,13).Value = 550
$ws.Cells.Item(435,14).Value = "`$/unidad"
$ws.Cells.Item(435,15).Value = "Región Metropolitana"
$ws.Cells.Item(435,16).Value = 550
$ws.Cells.Item(435,17).Value = 1
$ws.Cells.Item(435,18).Value = "Hortaliza"

# New row 436: Segunda quality for the new week (2021-11-16 = serial 44516)
$ws.Cells.Item(436,1).Value = 9
$ws.Cells.Item(436,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(436,3).Value = "Metropolitana"
$ws.Cells.Item(436,4).Value = 44516
$ws.Cells.Item(436,5).Value = 13
$ws.Cells.Item(436,6).Value = 100112023
$ws.Cells.Item(436,7).Value = "Brócoli"
$ws.Cells.Item(436,8).Value = "Sin especificar"
$ws.Cells.Item(436,9).Value = "Segunda"
$ws.Cells.Item(436,10).Value = 2500
$ws.Cells.Item(436,11).Value = 400
$ws.Cells.Item(436,12).Value = 400
$ws.Cells.Item(436,13).Value = 400
$ws.Cells.Item(436,14).Value = "`$/unidad"
$ws.Cells.Item(436,15).Value = "Región Metropolitana"
$ws.Cells.Item(436,16).Value = 400
$ws.Cells.Item(436,17).Value = 1
$ws.Cells.Item(436,18).Value = "Hortaliza"
